$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H51").Value = 16789.1
$ws.Range("I51").Value = 9787
$ws.Range("J51").Value = 21457.166
$ws.Range("K51").Value = 9787
$ws.Range("L51").Value = 21457.166
$ws.Range("M51").Value = -9303
$ws.Range("N51").Value = -22425.166
$ws.Range("H107").Value = 2069.0667
$ws.Range("I107").Value = 1894.909
$ws.Range("J107").Value = 2548
$ws.Range("K107").Value = 1894.909
$ws.Range("L107").Value = 2548
$ws.Range("M107").Value = 25.09099999999989
$ws.Range("N107").Value = -6388
$ws.Range("H111").Value = 4481.6665
$ws.Range("I111").Value = 3518.75
$ws.Range("K111").Value = 10556.25
$ws.Range("M111").Value = -7489.25
$ws.Range("H116").Value = 4999.5
$ws.Range("J116").Value = 4999
$ws.Range("L116").Value = 4999
$ws.Range("N116").Value = -11883
$ws.Range("H131").Value = 680.6667
$ws.Range("I131").Value = 495
$ws.Range("K131").Value = 1485
$ws.Range("M131").Value = 3555
$ws.Range("H132").Value = 2399.372
$ws.Range("I132").Value = 2363.8157
$ws.Range("K132").Value = 7091.4471
$ws.Range("M132").Value = -4561.4471

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H4").Value = 625
$ws.Range("J4").Value = 887.5
$ws.Range("L4").Value = 887.5
$ws.Range("N4").Value = -1119.5
$ws.Range("H32").Value = 9795.846
$ws.Range("I32").Value = 784.4483
$ws.Range("K32").Value = 784.4483
$ws.Range("M32").Value = -497.4483
$ws.Range("H45").Value = 3337.9
$ws.Range("I45").Value = 2825.7273
$ws.Range("K45").Value = 2825.7273
$ws.Range("M45").Value = -2448.7273
$ws.Range("H88").Value = 1998.8889
$ws.Range("I88").Value = 1648
$ws.Range("J88").Value = 2174.3333
$ws.Range("K88").Value = 1648
$ws.Range("L88").Value = 2174.3333
$ws.Range("M88").Value = -1242
$ws.Range("N88").Value = -2986.3333
$ws.Range("H91").Value = 1998.8889
$ws.Range("I91").Value = 1648
$ws.Range("J91").Value = 2174.3333
$ws.Range("K91").Value = 1648
$ws.Range("L91").Value = 2174.3333
$ws.Range("M91").Value = -244
$ws.Range("N91").Value = -4982.3333
$ws.Range("H97").Value = 1891.9231
$ws.Range("I97").Value = 585.13336
$ws.Range("J97").Value = 3673.9092
$ws.Range("K97").Value = 585.13336
$ws.Range("L97").Value = 3673.9092
$ws.Range("M97").Value = -89.13336000000004
$ws.Range("N97").Value = -4665.9092
$ws.Range("H122").Value = 3260.611
$ws.Range("I122").Value = 3089.5386
$ws.Range("J122").Value = 3705.4
$ws.Range("K122").Value = 9268.6158
$ws.Range("L122").Value = 11116.2
$ws.Range("M122").Value = -6818.6158
$ws.Range("N122").Value = -16016.2
$ws.Range("H125").Value = 54500
$ws.Range("J125").Value = 54500
$ws.Range("L125").Value = 54500
$ws.Range("N125").Value = -64340

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H22").Value = 368.33334
$ws.Range("I22").Value = 275
$ws.Range("J22").Value = 555
$ws.Range("K22").Value = 275
$ws.Range("L22").Value = 555
$ws.Range("M22").Value = -102
$ws.Range("N22").Value = -901
$ws.Range("H76").Value = 53333
$ws.Range("J76").Value = 53333
$ws.Range("L76").Value = 53333
$ws.Range("N76").Value = -53963
$ws.Range("H79").Value = 53333
$ws.Range("J79").Value = 53333
$ws.Range("L79").Value = 53333
$ws.Range("N79").Value = -55517
$ws.Range("H105").Value = 3491.818
$ws.Range("I105").Value = 3062.5
$ws.Range("J105").Value = 4636.6665
$ws.Range("K105").Value = 3062.5
$ws.Range("L105").Value = 4636.6665
$ws.Range("M105").Value = -1315.5
$ws.Range("N105").Value = -8130.6665

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H62").Value = 5533.8335
$ws.Range("J62").Value = 7248.75
$ws.Range("L62").Value = 7248.75
$ws.Range("N62").Value = -8496.75
$ws.Range("H65").Value = 5533.8335
$ws.Range("J65").Value = 7248.75
$ws.Range("L65").Value = 36243.75
$ws.Range("N65").Value = -42483.75
$ws.Range("H107").Value = 1716.0625
$ws.Range("I107").Value = 1753.6666
$ws.Range("J107").Value = 1603.25
$ws.Range("K107").Value = 1753.6666
$ws.Range("L107").Value = 1603.25
$ws.Range("M107").Value = 166.3334
$ws.Range("N107").Value = -5443.25
$ws.Range("H111").Value = 1000
$ws.Range("I111").Value = 1000
$ws.Range("K111").Value = 1000
$ws.Range("M111").Value = 3090

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H46").Value = 92985.82000000001
$ws.Range("I46").Value = 2559.8
$ws.Range("K46").Value = 7679.400000000001
$ws.Range("M46").Value = -7588.400000000001
$ws.Range("H70").Value = 8832.714
$ws.Range("I70").Value = 8832.714
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 26498.142
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -26183.142
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 8832.714
$ws.Range("I73").Value = 8832.714
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 26498.142
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -25406.142
$ws.Range("N73").ClearContents()
$ws.Range("H80").Value = 3999.5
$ws.Range("I80").Value = 5000
$ws.Range("K80").Value = 15000
$ws.Range("M80").Value = -14064
$ws.Range("H83").Value = 3999.5
$ws.Range("I83").Value = 5000
$ws.Range("K83").Value = 45000
$ws.Range("M83").Value = -40320
$ws.Range("H86").Value = 33798.668
$ws.Range("I86").Value = 698
$ws.Range("K86").Value = 2094
$ws.Range("M86").Value = -908
$ws.Range("H89").Value = 33798.668
$ws.Range("I89").Value = 698
$ws.Range("K89").Value = 6282
$ws.Range("M89").Value = -354
$ws.Range("H113").Value = 3127.8572
$ws.Range("I113").Value = 2099.3333
$ws.Range("K113").Value = 6297.999899999999
$ws.Range("M113").Value = -4127.999899999999
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("M119").ClearContents()
$ws.Range("H122").Value = 33999
$ws.Range("J122").Value = 33999
$ws.Range("L122").Value = 305991
$ws.Range("N122").Value = -310891
$ws.Range("H124").Value = 2746
$ws.Range("I124").Value = 2746
$ws.Range("K124").Value = 8238
$ws.Range("M124").Value = -3328

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H132").Value = 2902.56
$ws.Range("I132").Value = 3080.2104
$ws.Range("J132").Value = 2340
$ws.Range("K132").Value = 9240.6312
$ws.Range("L132").Value = 7020
$ws.Range("M132").Value = -6710.6312
$ws.Range("N132").Value = -12080

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H22").Value = 2191.5334
$ws.Range("I22").Value = 2033.7333
$ws.Range("J22").Value = 2349.3333
$ws.Range("K22").Value = 2033.7333
$ws.Range("L22").Value = 2349.3333
$ws.Range("M22").Value = -1738.7333
$ws.Range("N22").Value = -2939.3333
$ws.Range("H27").Value = 2191.5334
$ws.Range("I27").Value = 2033.7333
$ws.Range("J27").Value = 2349.3333
$ws.Range("K27").Value = 2033.7333
$ws.Range("L27").Value = 2349.3333
$ws.Range("M27").Value = -1926.7333
$ws.Range("N27").Value = -2563.3333
$ws.Range("H40").Value = 8179.2
$ws.Range("I40").Value = 5151.9443
$ws.Range("K40").Value = 5151.9443
$ws.Range("M40").Value = -5015.9443
$ws.Range("H100").Value = 4117.143
$ws.Range("I100").Value = 2591.5715
$ws.Range("J100").Value = 5642.7144
$ws.Range("K100").Value = 2591.5715
$ws.Range("L100").Value = 5642.7144
$ws.Range("M100").Value = -2050.5715
$ws.Range("N100").Value = -6724.7144
$ws.Range("H122").Value = 7360.615
$ws.Range("J122").Value = 11823.833
$ws.Range("L122").Value = 35471.499
$ws.Range("N122").Value = -40371.499
$ws.Range("H136").Value = 1420.6364
$ws.Range("I136").Value = 1302.875
$ws.Range("J136").Value = 1734.6666
$ws.Range("K136").Value = 3908.625
$ws.Range("L136").Value = 5203.9998
$ws.Range("M136").Value = -1358.625
$ws.Range("N136").Value = -10303.9998

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H62").Value = 10000
$ws.Range("J62").Value = 10000
$ws.Range("L62").Value = 10000
$ws.Range("N62").Value = -11248
$ws.Range("H65").Value = 10000
$ws.Range("J65").Value = 10000
$ws.Range("L65").Value = 50000
$ws.Range("N65").Value = -56240
$ws.Range("H81").Value = 4346.05
$ws.Range("J81").Value = 5000
$ws.Range("L81").Value = 10000
$ws.Range("N81").Value = -12122
$ws.Range("H84").Value = 4346.05
$ws.Range("J84").Value = 5000
$ws.Range("L84").Value = 50000
$ws.Range("N84").Value = -60608
$ws.Range("H126").Value = 2376.818
$ws.Range("I126").Value = 2375
$ws.Range("K126").Value = 7125
$ws.Range("M126").Value = -4655
$ws.Range("H136").Value = 33799.8
$ws.Range("I136").Value = 32249.75
$ws.Range("K136").Value = 96749.25
$ws.Range("M136").Value = -94199.25
